# Session added for debugging bug #98.
# Update the LB/UB acceptance-criteria bounds on Sheet1 (rows 2-5, columns D/E)
# so every row now uses the "no dose" bounds of 0 / 1, and move the active
# selection to E6 as left by the editor.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# LB (column D) -> 0, UB (column E) -> 1 for data rows 2 through 5
$ws1.Range("D2:D5").Value = 0
$ws1.Range("E2:E5").Value = 1

# Leave the selection where the author left it when saving the session.
$ws1.Activate()
$ws1.Range("E6").Select()
